# "new code for discount and impersonate a customer"
#
# 1. On the existing "Discount" sheet, the Discountable-Item column (L1)
#    is changed from "Plan Product" to "plan", and the order-item column
#    (M1) is changed from "Product Code1" to "Product Code1 Description".
# 2. A brand-new "ImpersonateCust" worksheet is inserted right after the
#    "Discount" sheet (and therefore right before "AddingAsset"), holding
#    the login/company/username data needed to impersonate a customer.

$wb = $excel.ActiveWorkbook

# --- 1. Update the Discount sheet ---------------------------------------
$discountSheet = $wb.Worksheets.Item("Discount")
$discountSheet.Range("L1").Value = "plan"
$discountSheet.Range("M1").Value = "Product Code1 Description"
$discountSheet.Range("M1").Select()

# --- 2. Insert the new ImpersonateCust sheet -----------------------------
$newSheet = $wb.Worksheets.Add($null, $discountSheet)
$newSheet.Name = "ImpersonateCust"
$newSheet.Range("A1").Value = "admin"
$newSheet.Range("B1").Value = "Webdata@123"
$newSheet.Range("C1").Value = "Web Data US"
$newSheet.Range("D1").Value = "James"
$newSheet.Range("D1").Select()

# Keep "Discount" as the active/tab-selected sheet, as it was before.
$discountSheet.Activate()
$discountSheet.Range("M1").Select()
